# The document stores explicit "\n" markers as literal text runs (not real
# line breaks) at the end of most paragraphs, with blank spacer paragraphs
# (containing only "\n") in between the real content paragraphs.
#
# Target change: remove the blank spacer paragraph and the whole
# "Please read Using Pressure Canners..." paragraph that come right after
# the "Quality: Select mature, dry seeds..." paragraph, so that paragraph
# flows straight into "Procedure: Place dried beans or peas...". The
# "Quality..." paragraph itself is left completely untouched.

$d = $word.ActiveDocument

# Carriage return (Chr 13) represents the actual paragraph mark between
# paragraphs; it is distinct from the literal "\n" text runs.
$cr = [char]13

$removedParaText = "Please read Using Pressure Canners before beginning. If this is your first time canning, it is recommended that you read Principles of Home Canning."

# This spans: blank paragraph's "\n" run -> its paragraph mark -> the
# "Please read..." text -> its own "\n" run -> its paragraph mark. It
# starts right after the paragraph mark that already ends the
# "Quality..." paragraph, so that paragraph mark (and everything before
# it) is left untouched. Deleting this match merges what is left on
# either side, i.e. the "Quality..." paragraph (untouched) directly
# followed by the paragraph that used to contain "Procedure...".
$searchText = "\n" + $cr + $removedParaText + "\n" + $cr

$find = $d.Content.Find
$find.ClearFormatting()
$result = $find.Execute($searchText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

if ($result) {
    Write-Host "Removed blank paragraph and 'Please read Using Pressure Canners...' paragraph."
} else {
    Write-Host "WARNING: target text not found; document left unchanged."
}
